$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 53, shifting existing rows 53:69 down to 54:70
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly price record
$ws.Range("A53").Value = 2
$ws.Range("B53").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C53").Value = "Coquimbo"
$ws.Range("D53").Value = 44798
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = 100112022
$ws.Range("G53").Value = "Arveja Verde"
$ws.Range("H53").Value = "Perfection"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 400
$ws.Range("K53").Value = 30000
$ws.Range("L53").Value = 32000
$ws.Range("M53").Value = 31000
$ws.Range("N53").Value = "$/malla 25 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 1240
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"
